$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet contains a structured table ("Tableau1") covering B4:L44.
# Add a new row to the table; this automatically grows the table/autoFilter
# range (B4:L44 -> B4:L45) and the sheet dimension.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Copy the date formatting from the row above (B44) so the new date cell
# (B45) reuses the existing date style instead of creating a new one.
$ws.Range("B44").Copy()
$ws.Range("B45").PasteSpecial(-4122)

# Fill in the new journal entry (row 45).
$ws.Range("B45").Value = 44267
$ws.Range("C45").Value = 0.61458333333333337
$ws.Range("D45").Value = 0.625
$ws.Range("E45").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F45").Value = "Ma-20"
$ws.Range("G45").Value = "Code"
$ws.Range("H45").Value = "Menu"
$ws.Range("I45").Value = "CPNV"
$ws.Range("J45").Value = "J'ai temriner mon menu, on ne peut pas sortir du programme si on ne choisis pas 0"
$ws.Range("K45").Value = "Oui"

# Match the wrapped-text row height used by similar long "Descriptif" rows.
$ws.Rows.Item(45).RowHeight = 43.2

# Move the selection to the newly added row, like the author's last edit.
$ws.Range("L45").Select()
